$wb = $excel.ActiveWorkbook

# --- BRPSPTY data sheet: zero out the historical (2016-2019) RPS values ---
$data = $wb.Worksheets.Item("BRPSPTY")
$data.Range("B2").Value = 0
$data.Range("C2").Value = 0
$data.Range("D2").Value = 0
$data.Range("E2").Value = 0

# --- About sheet: append explanatory note rows ---
$about = $wb.Worksheets.Item("About")
$about.Range("A21").Value = "RPO national targets have historically not been achieved for wind and solar generation."
$about.Range("A22").Value = "The partial target achieved in 2018 & 2019, is accounted in the real-world capacity deployed for  "
$about.Range("A23").Value = "wind and solar in elec/BPMCCS. Hence, we set the historical targets till 2019 to zero in this variable."
